$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 13.20963693679378
$ws.Range("D2").Value = 7.881186435122017
$ws.Range("E2").Value = 13.6945360507812
$ws.Range("F2").Value = 40.69777765117924
$ws.Range("G2").Value = 3.70654522015496
$ws.Range("J2").Value = 10.80080169445698
$ws.Range("K2").Value = 17.35962654059482
$ws.Range("L2").Value = 9.46195815782611
$ws.Range("M2").Value = 19.78728721444257
$ws.Range("N2").Value = 19.64231612616365
$ws.Range("O2").Value = 31.43149510503807

$ws.Range("C3").Value = 13.19284855831593
$ws.Range("D3").Value = 7.870321294100675
$ws.Range("E3").Value = 13.7164388702485
$ws.Range("F3").Value = 40.78182227029627
$ws.Range("G3").Value = 3.708886641905894
$ws.Range("J3").Value = 10.8251291443859
$ws.Range("K3").Value = 17.07190179699658
$ws.Range("L3").Value = 9.475756496461733
$ws.Range("M3").Value = 19.68642967762725
$ws.Range("N3").Value = 19.69714928907865
$ws.Range("O3").Value = 31.51899584697338

$ws.Range("C4").Value = 13.18509345778682
$ws.Range("D4").Value = 7.864535925870994
$ws.Range("E4").Value = 13.73156341128211
$ws.Range("F4").Value = 40.84229094211128
$ws.Range("G4").Value = 3.71040108481766
$ws.Range("J4").Value = 10.84100278575613
$ws.Range("K4").Value = 16.89527845037528
$ws.Range("L4").Value = 9.484739571326047
$ws.Range("M4").Value = 19.62701506254915
$ws.Range("N4").Value = 19.73268595868009
$ws.Range("O4").Value = 31.57915267924656

$ws.Range("C5").Value = 13.18257799541644
$ws.Range("D5").Value = 7.862402845898749
$ws.Range("E5").Value = 13.73814864460826
$ws.Range("F5").Value = 40.86915712363207
$ws.Range("G5").Value = 3.711037606983839
$ws.Range("J5").Value = 10.84770741606979
$ws.Range("K5").Value = 16.82340100110103
$ws.Range("L5").Value = 9.488529064689233
$ws.Range("M5").Value = 19.60345247203895
$ws.Range("N5").Value = 19.74763850268974
$ws.Range("O5").Value = 31.60528095185629

$ws.Range("C6").Value = 13.1821993232159
$ws.Range("D6").Value = 7.862062257988418
$ws.Range("E6").Value = 13.73926760764684
$ws.Range("F6").Value = 40.87375245988213
$ws.Range("G6").Value = 3.711144473001485
$ws.Range("J6").Value = 10.8488349830907
$ws.Range("K6").Value = 16.81147430639786
$ws.Range("L6").Value = 9.489166098111774
$ws.Range("M6").Value = 19.59957963925804
$ws.Range("N6").Value = 19.75014984533354
$ws.Range("O6").Value = 31.6097169072676

$ws.Range("C7").Value = 13.18505691918646
$ws.Range("D7").Value = 7.86450624716548
$ws.Range("E7").Value = 13.73165051345283
$ws.Range("F7").Value = 40.84264426735887
$ws.Range("G7").Value = 3.71040959064923
$ws.Range("J7").Value = 10.84109225054927
$ws.Range("K7").Value = 16.8943085725756
$ws.Range("L7").Value = 9.484790155698636
$ws.Range("M7").Value = 19.62669463739339
$ws.Range("N7").Value = 19.73288570512218
$ws.Range("O7").Value = 31.57949852421233

$ws.Range("C8").Value = 13.20332024470788
$ws.Range("D8").Value = 7.877257348159977
$ws.Range("E8").Value = 13.70174050288348
$ws.Range("F8").Value = 40.72491377796565
$ws.Range("G8").Value = 3.707336639975736
$ws.Range("J8").Value = 10.80899573465157
$ws.Range("K8").Value = 17.26046287820109
$ws.Range("L8").Value = 9.466610018634803
$ws.Range("M8").Value = 19.75200090175288
$ws.Range("N8").Value = 19.66083538345036
$ws.Range("O8").Value = 31.46032869070396

$ws.Range("C9").Value = 13.25924971663992
$ws.Range("D9").Value = 7.909211036473248
$ws.Range("E9").Value = 13.65636895480938
$ws.Range("F9").Value = 40.56456558366426
$ws.Range("G9").Value = 3.701917130667086
$ws.Range("J9").Value = 10.75346237977541
$ws.Range("K9").Value = 17.97502642532358
$ws.Range("L9").Value = 9.434995716212743
$ws.Range("M9").Value = 20.0168530712047
$ws.Range("N9").Value = 19.53432313601002
$ws.Range("O9").Value = 31.27780453372159

$ws.Range("C10").Value = 13.31238740212929
$ws.Range("D10").Value = 7.936810972478272
$ws.Range("E10").Value = 13.63110742164126
$ws.Range("F10").Value = 40.48997118991144
$ws.Range("G10").Value = 3.698301206287086
$ws.Range("J10").Value = 10.71714644781798
$ws.Range("K10").Value = 18.4929687301218
$ws.Range("L10").Value = 9.414206980700369
$ws.Range("M10").Value = 20.22197428672398
$ws.Range("N10").Value = 19.45031355753176
$ws.Range("O10").Value = 31.1750806257956

$ws.Range("C11").Value = 13.33912161509298
$ws.Range("D11").Value = 7.950236157927526
$ws.Range("E11").Value = 13.62136257163345
$ws.Range("F11").Value = 40.46545616682551
$ws.Range("G11").Value = 3.696734811288499
$ws.Range("J11").Value = 10.70159249244172
$ws.Range("K11").Value = 18.7260696750887
$ws.Range("L11").Value = 9.405274290966682
$ws.Range("M11").Value = 20.31732415688667
$ws.Range("N11").Value = 19.41402179327268
$ws.Range("O11").Value = 31.13519754841904

$ws.Range("C12").Value = 13.34960811586207
$ws.Range("D12").Value = 7.955442599264165
$ws.Range("E12").Value = 13.61792309377881
$ws.Range("F12").Value = 40.45752925551705
$ws.Range("G12").Value = 3.696152883326093
$ws.Range("J12").Value = 10.69584107490874
$ws.Range("K12").Value = 18.8138988511379
$ws.Range("L12").Value = 9.401966727144726
$ws.Range("M12").Value = 20.35370179438007
$ws.Range("N12").Value = 19.4005547374367
$ws.Range("O12").Value = 31.12108196313059

$ws.Range("C13").Value = 13.34733361558141
$ws.Range("D13").Value = 7.954315885991556
$ws.Range("E13").Value = 13.61865270585277
$ws.Range("F13").Value = 40.45917610564456
$ws.Range("G13").Value = 3.696277713295416
$ws.Range("J13").Value = 10.69707359132527
$ws.Range("K13").Value = 18.79500411292172
$ws.Range("L13").Value = 9.402675737446389
$ws.Range("M13").Value = 20.34585554208939
$ws.Range("N13").Value = 19.40344285470668
$ws.Range("O13").Value = 31.12407805660665

$ws.Range("C14").Value = 13.33997710721931
$ws.Range("D14").Value = 7.950662052092411
$ws.Range("E14").Value = 13.62107458257851
$ws.Range("F14").Value = 40.46477682135287
$ws.Range("G14").Value = 3.696686710949609
$ws.Range("J14").Value = 10.70111654600612
$ws.Range("K14").Value = 18.7333047174227
$ws.Range("L14").Value = 9.405000673421837
$ws.Range("M14").Value = 20.32031167197025
$ws.Range("N14").Value = 19.41290832721536
$ws.Range("O14").Value = 31.13401644944959

$ws.Range("C15").Value = 13.33551811403863
$ws.Range("D15").Value = 7.948439866473117
$ws.Range("E15").Value = 13.62259068331779
$ws.Range("F15").Value = 40.4683841143878
$ws.Range("G15").Value = 3.696938694841525
$ws.Range("J15").Value = 10.70361100019431
$ws.Range("K15").Value = 18.69545231227378
$ws.Range("L15").Value = 9.406434528101096
$ws.Range("M15").Value = 20.30469989330123
$ws.Range("N15").Value = 19.41874209933876
$ws.Range("O15").Value = 31.1402326513296

$ws.Range("C16").Value = 13.3106913423569
$ws.Range("D16").Value = 7.935950899876857
$ws.Range("E16").Value = 13.63177937808123
$ws.Range("F16").Value = 40.49176301488083
$ws.Range("G16").Value = 3.698405148731462
$ws.Range("J16").Value = 10.71818234369291
$ws.Range("K16").Value = 18.47767756369105
$ws.Range("L16").Value = 9.414801273253353
$ws.Range("M16").Value = 20.21578206634896
$ws.Range("N16").Value = 19.45272395956431
$ws.Range("O16").Value = 31.17782508864963

$ws.Range("C17").Value = 13.29611321997203
$ws.Range("D17").Value = 7.928510387231495
$ws.Range("E17").Value = 13.63786338833766
$ws.Range("F17").Value = 40.50851913193257
$ws.Range("G17").Value = 3.69932483727287
$ws.Range("J17").Value = 10.72736858614957
$ws.Range("K17").Value = 18.34337819963264
$ws.Range("L17").Value = 9.420068028310819
$ws.Range("M17").Value = 20.16174025522437
$ws.Range("N17").Value = 19.47406301565646
$ws.Range("O17").Value = 31.20264242199855

$ws.Range("C18").Value = 13.28796973881888
$ws.Range("D18").Value = 7.924312828471362
$ws.Range("E18").Value = 13.64152719423973
$ws.Range("F18").Value = 40.51904322236257
$ws.Range("G18").Value = 3.699861210663121
$ws.Range("J18").Value = 10.73274325214302
$ws.Range("K18").Value = 18.26590070545619
$ws.Range("L18").Value = 9.423146687155057
$ws.Range("M18").Value = 20.13085015043798
$ws.Range("N18").Value = 19.48651788302217
$ws.Range("O18").Value = 31.21756097479534

$ws.Range("C19").Value = 13.28525412983401
$ws.Range("D19").Value = 7.922905768925172
$ws.Range("E19").Value = 13.64279595387711
$ws.Range("F19").Value = 40.52275866860359
$ws.Range("G19").Value = 3.70004408896556
$ws.Range("J19").Value = 10.73457866123824
$ws.Range("K19").Value = 18.239630788285
$ws.Range("L19").Value = 9.424197555947059
$ws.Range("M19").Value = 20.12042517338076
$ws.Range("N19").Value = 19.49076603763769
$ws.Range("O19").Value = 31.22272270859846

$ws.Range("C20").Value = 13.29764013347682
$ws.Range("D20").Value = 7.929293970532093
$ws.Range("E20").Value = 13.63719871918124
$ws.Range("F20").Value = 40.50664365716495
$ws.Range("G20").Value = 3.699226170140511
$ws.Range("J20").Value = 10.7263812817723
$ws.Range("K20").Value = 18.35769920122808
$ws.Range("L20").Value = 9.419502266995071
$ws.Range("M20").Value = 20.16747325268956
$ws.Range("N20").Value = 19.47177268914165
$ws.Range("O20").Value = 31.19993387337505

$ws.Range("C21").Value = 13.34212809094343
$ws.Range("D21").Value = 7.951731964105943
$ws.Range("E21").Value = 13.62035641927649
$ws.Range("F21").Value = 40.46309492995596
$ws.Range("G21").Value = 3.69656627405155
$ws.Range("J21").Value = 10.69992527601739
$ws.Range("K21").Value = 18.75143991295414
$ws.Range("L21").Value = 9.404315749087838
$ws.Range("M21").Value = 20.32780736024024
$ws.Range("N21").Value = 19.41012060920202
$ws.Range("O21").Value = 31.13107048758985

$ws.Range("C22").Value = 13.37331556172422
$ws.Range("D22").Value = 7.967109968602736
$ws.Range("E22").Value = 13.61080993285545
$ws.Range("F22").Value = 40.44253993837041
$ws.Range("G22").Value = 3.694893319560209
$ws.Range("J22").Value = 10.68344199694426
$ws.Range("K22").Value = 19.00616692075477
$ws.Range("L22").Value = 9.394827791164177
$ws.Range("M22").Value = 20.43416217276868
$ws.Range("N22").Value = 19.37143485634822
$ws.Range("O22").Value = 31.09181996022957

$ws.Range("C23").Value = 13.35647890523958
$ws.Range("D23").Value = 7.958837995823885
$ws.Range("E23").Value = 13.61577157020902
$ws.Range("F23").Value = 40.45278655332982
$ws.Range("G23").Value = 3.695780237819493
$ws.Range("J23").Value = 10.69216571075799
$ws.Range("K23").Value = 18.87047805190235
$ws.Range("L23").Value = 9.399851786944684
$ws.Range("M23").Value = 20.37726274661194
$ws.Range("N23").Value = 19.39193538396505
$ws.Range("O23").Value = 31.11224121398741

$ws.Range("C24").Value = 13.29694907514682
$ws.Range("D24").Value = 7.928939462874612
$ws.Range("E24").Value = 13.63749869891486
$ws.Range("F24").Value = 40.5074887847118
$ws.Range("G24").Value = 3.699270753771821
$ws.Range("J24").Value = 10.72682735119571
$ws.Range("K24").Value = 18.35122550192735
$ws.Range("L24").Value = 9.419757889637918
$ws.Range("M24").Value = 20.16488080445356
$ws.Range("N24").Value = 19.47280756392377
$ws.Range("O24").Value = 31.20115638136269

$ws.Range("C25").Value = 13.24198640533599
$ws.Range("D25").Value = 7.899833988377079
$ws.Range("E25").Value = 13.66722359765846
$ws.Range("F25").Value = 40.60036917400154
$ws.Range("G25").Value = 3.703318729775767
$ws.Range("J25").Value = 10.76769585219535
$ws.Range("K25").Value = 17.78260795448779
$ws.Range("L25").Value = 9.443118403444533
$ws.Range("M25").Value = 19.53432313601002
$ws.Range("N25").Value = 19.56697328540271
$ws.Range("O25").Value = 31.32168460738745
